$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells H1/I1 with values matching the added columns
$ws.Range("H1").Value = "mdescription"
$ws.Range("I1").Value = "IP"

# Match the bold/centered/bordered header style already used by A1:G1
# (copy the formatting from the existing header cell G1).
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# New data row 4 (Fanuc / EU 63 machine with its mdescription + IP)
$ws.Range("B4").Value = "Fanuc"
$ws.Range("C4").Value = "EU 63"
$ws.Range("H4").Value = "M614"
$ws.Range("I4").Value = "x.x.x.x.x"
